# Update "SOPORTE DE ESTANTE" price list
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update date in A1 (one month later: 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = Get-Date -Year 2024 -Month 5 -Day 24 -Hour 0 -Minute 0 -Second 0

# Updated prices (NEGRO section, rows 26-32)
$ws.Range("D26").Value = 6797.942
$ws.Range("D27").Value = 8507.816999999999
$ws.Range("D28").Value = 11071.228
$ws.Range("D29").Value = 15643.851
$ws.Range("D30").Value = 23056.99
$ws.Range("D31").Value = 31703.345
$ws.Range("D32").Value = 40757.109

# Updated prices (BLANCO section, rows 34-40)
$ws.Range("D34").Value = 6797.942
$ws.Range("D35").Value = 8507.816999999999
$ws.Range("D36").Value = 11071.228
$ws.Range("D37").Value = 15643.851
$ws.Range("D38").Value = 23056.99
$ws.Range("D39").Value = 31703.345
$ws.Range("D40").Value = 40757.109
